$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 15,20
$data[0,0] = "ECs"
$data[0,1] = "Fasl"
$data[0,2] = "Fas"
$data[0,3] = "ECs"
$data[0,4] = 1
$data[0,5] = 0.3333333333333333
$data[0,6] = 0.068049
$data[0,7] = 0.204147
$data[0,8] = 0.1675230998998868
$data[0,9] = 0.1675230998998868
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 3.042047666666666
$data[0,13] = 9.126142999999999
$data[0,14] = 0.1400448830218481
$data[0,15] = 0.140044883021848
$data[0,16] = 0.207008301669
$data[0,17] = 1.863074715021
$data[0,18] = 0.02346075292893702
$data[0,19] = 0.023460752928937

$data[1,0] = "ECs"
$data[1,1] = "Fasl"
$data[1,2] = "Fas"
$data[1,3] = "FAPs"
$data[1,4] = 1
$data[1,5] = 0.3333333333333333
$data[1,6] = 0.068049
$data[1,7] = 0.204147
$data[1,8] = 0.1675230998998868
$data[1,9] = 0.1675230998998868
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 9.854944999999999
$data[1,13] = 29.564835
$data[1,14] = 0.4536860598322029
$data[1,15] = 0.4536860598322028
$data[1,16] = 0.6706191523049999
$data[1,17] = 6.035572370744999
$data[1,18] = 0.07600289512445614
$data[1,19] = 0.07600289512445611

$data[2,0] = "ECs"
$data[2,1] = "Fasl"
$data[2,2] = "Fas"
$data[2,3] = "Inflammatory-Mac"
$data[2,4] = 1
$data[2,5] = 0.3333333333333333
$data[2,6] = 0.068049
$data[2,7] = 0.204147
$data[2,8] = 0.1675230998998868
$data[2,9] = 0.1675230998998868
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 4.111990666666666
$data[2,13] = 12.335972
$data[2,14] = 0.1893011928150581
$data[2,15] = 0.1893011928150581
$data[2,16] = 0.2798168528759999
$data[2,17] = 2.518351675884
$data[2,18] = 0.03171232263512472
$data[2,19] = 0.0317123226351247

$data[3,0] = "ECs"
$data[3,1] = "Fasl"
$data[3,2] = "Fas"
$data[3,3] = "MuSCs"
$data[3,4] = 1
$data[3,5] = 0.3333333333333333
$data[3,6] = 0.068049
$data[3,7] = 0.204147
$data[3,8] = 0.1675230998998868
$data[3,9] = 0.1675230998998868
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 2.157352
$data[3,13] = 6.472056
$data[3,14] = 0.09931669111812624
$data[3,15] = 0.09931669111812623
$data[3,16] = 0.146805646248
$data[3,17] = 1.321250816232
$data[3,18] = 0.01663783996790806
$data[3,19] = 0.01663783996790806

$data[4,0] = "ECs"
$data[4,1] = "Fasl"
$data[4,2] = "Fas"
$data[4,3] = "Resolving-Mac"
$data[4,4] = 1
$data[4,5] = 0.3333333333333333
$data[4,6] = 0.068049
$data[4,7] = 0.204147
$data[4,8] = 0.1675230998998868
$data[4,9] = 0.1675230998998868
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 2.555612666666667
$data[4,13] = 7.666838
$data[4,14] = 0.1176511732127647
$data[4,15] = 0.1176511732127646
$data[4,16] = 0.173906886354
$data[4,17] = 1.565161977186
$data[4,18] = 0.01970928924346086
$data[4,19] = 0.01970928924346085

$data[5,0] = "FAPs"
$data[5,1] = "Fasl"
$data[5,2] = "Fas"
$data[5,3] = "ECs"
$data[5,4] = 1
$data[5,5] = 0.3333333333333333
$data[5,6] = 0.298478
$data[5,7] = 0.895434
$data[5,8] = 0.7347934548915988
$data[5,9] = 0.7347934548915986
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 3.042047666666666
$data[5,13] = 9.126142999999999
$data[5,14] = 0.1400448830218481
$data[5,15] = 0.140044883021848
$data[5,16] = 0.9079843034513332
$data[5,17] = 8.171858731061999
$data[5,18] = 0.1029040634355136
$data[5,19] = 0.1029040634355135

$data[6,0] = "FAPs"
$data[6,1] = "Fasl"
$data[6,2] = "Fas"
$data[6,3] = "FAPs"
$data[6,4] = 1
$data[6,5] = 0.3333333333333333
$data[6,6] = 0.298478
$data[6,7] = 0.895434
$data[6,8] = 0.7347934548915988
$data[6,9] = 0.7347934548915986
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 9.854944999999999
$data[6,13] = 29.564835
$data[6,14] = 0.4536860598322029
$data[6,15] = 0.4536860598322028
$data[6,16] = 2.94148427371
$data[6,17] = 26.47335846339
$data[6,18] = 0.333365547340261
$data[6,19] = 0.3333655473402609

$data[7,0] = "FAPs"
$data[7,1] = "Fasl"
$data[7,2] = "Fas"
$data[7,3] = "Inflammatory-Mac"
$data[7,4] = 1
$data[7,5] = 0.3333333333333333
$data[7,6] = 0.298478
$data[7,7] = 0.895434
$data[7,8] = 0.7347934548915988
$data[7,9] = 0.7347934548915986
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 4.111990666666666
$data[7,13] = 12.335972
$data[7,14] = 0.1893011928150581
$data[7,15] = 0.1893011928150581
$data[7,16] = 1.227338750205333
$data[7,17] = 11.046048751848
$data[7,18] = 0.1390972774836773
$data[7,19] = 0.1390972774836772

$data[8,0] = "FAPs"
$data[8,1] = "Fasl"
$data[8,2] = "Fas"
$data[8,3] = "MuSCs"
$data[8,4] = 1
$data[8,5] = 0.3333333333333333
$data[8,6] = 0.298478
$data[8,7] = 0.895434
$data[8,8] = 0.7347934548915988
$data[8,9] = 0.7347934548915986
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 2.157352
$data[8,13] = 6.472056
$data[8,14] = 0.09931669111812624
$data[8,15] = 0.09931669111812623
$data[8,16] = 0.6439221102559999
$data[8,17] = 5.795298992304
$data[8,18] = 0.07297725459508975
$data[8,19] = 0.07297725459508972

$data[9,0] = "FAPs"
$data[9,1] = "Fasl"
$data[9,2] = "Fas"
$data[9,3] = "Resolving-Mac"
$data[9,4] = 1
$data[9,5] = 0.3333333333333333
$data[9,6] = 0.298478
$data[9,7] = 0.895434
$data[9,8] = 0.7347934548915988
$data[9,9] = 0.7347934548915986
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 2.555612666666667
$data[9,13] = 7.666838
$data[9,14] = 0.1176511732127647
$data[9,15] = 0.1176511732127646
$data[9,16] = 0.7627941575213333
$data[9,17] = 6.865147417692
$data[9,18] = 0.08644931203705726
$data[9,19] = 0.08644931203705723

$data[10,0] = "MuSCs"
$data[10,1] = "Fasl"
$data[10,2] = "Fas"
$data[10,3] = "ECs"
$data[10,4] = 1
$data[10,5] = 0.3333333333333333
$data[10,6] = 0.03967966666666667
$data[10,7] = 0.119039
$data[10,8] = 0.09768344520851457
$data[10,9] = 0.09768344520851455
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 3.042047666666666
$data[10,13] = 9.126142999999999
$data[10,14] = 0.1400448830218481
$data[10,15] = 0.140044883021848
$data[10,16] = 0.1207074373974444
$data[10,17] = 1.086366936577
$data[10,18] = 0.01368006665739753
$data[10,19] = 0.01368006665739752

$data[11,0] = "MuSCs"
$data[11,1] = "Fasl"
$data[11,2] = "Fas"
$data[11,3] = "FAPs"
$data[11,4] = 1
$data[11,5] = 0.3333333333333333
$data[11,6] = 0.03967966666666667
$data[11,7] = 0.119039
$data[11,8] = 0.09768344520851457
$data[11,9] = 0.09768344520851455
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 9.854944999999999
$data[11,13] = 29.564835
$data[11,14] = 0.4536860598322029
$data[11,15] = 0.4536860598322028
$data[11,16] = 0.3910409326183333
$data[11,17] = 3.519368393565
$data[11,18] = 0.04431761736748586
$data[11,19] = 0.04431761736748584

$data[12,0] = "MuSCs"
$data[12,1] = "Fasl"
$data[12,2] = "Fas"
$data[12,3] = "Inflammatory-Mac"
$data[12,4] = 1
$data[12,5] = 0.3333333333333333
$data[12,6] = 0.03967966666666667
$data[12,7] = 0.119039
$data[12,8] = 0.09768344520851457
$data[12,9] = 0.09768344520851455
$data[12,10] = 3
$data[12,11] = 1
$data[12,12] = 4.111990666666666
$data[12,13] = 12.335972
$data[12,14] = 0.1893011928150581
$data[12,15] = 0.1893011928150581
$data[12,16] = 0.1631624189897778
$data[12,17] = 1.468461770908
$data[12,18] = 0.01849159269625619
$data[12,19] = 0.01849159269625618

$data[13,0] = "MuSCs"
$data[13,1] = "Fasl"
$data[13,2] = "Fas"
$data[13,3] = "MuSCs"
$data[13,4] = 1
$data[13,5] = 0.3333333333333333
$data[13,6] = 0.03967966666666667
$data[13,7] = 0.119039
$data[13,8] = 0.09768344520851457
$data[13,9] = 0.09768344520851455
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 2.157352
$data[13,13] = 6.472056
$data[13,14] = 0.09931669111812624
$data[13,15] = 0.09931669111812623
$data[13,16] = 0.08560300824266667
$data[13,17] = 0.770427074184
$data[13,18] = 0.009701596555128452
$data[13,19] = 0.009701596555128447

$data[14,0] = "MuSCs"
$data[14,1] = "Fasl"
$data[14,2] = "Fas"
$data[14,3] = "Resolving-Mac"
$data[14,4] = 1
$data[14,5] = 0.3333333333333333
$data[14,6] = 0.03967966666666667
$data[14,7] = 0.119039
$data[14,8] = 0.09768344520851457
$data[14,9] = 0.09768344520851455
$data[14,10] = 3
$data[14,11] = 1
$data[14,12] = 2.555612666666667
$data[14,13] = 7.666838
$data[14,14] = 0.1176511732127647
$data[14,15] = 0.1176511732127646
$data[14,16] = 0.1014058587424445
$data[14,17] = 0.9126527286820001
$data[14,18] = 0.01149257193224655
$data[14,19] = 0.01149257193224655

$rng = $ws.Range("A2:T16")
$rng.Value2 = $data

Write-Output "done"